# TC10_Canine_Filter_Breed-BostonTerr.xlsx
# "Fixed variables and query errors in Bread from TC01 to TC30"
#
# The Cases-tab query in cell B2 had a stray `Cohort` column/clause that
# needs to be dropped (the Cohort concept doesn't belong to this query),
# and the view needs to reflect the zoom level / selection that was active
# when the fix was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the Cases query in B2: remove the trailing Cohort column ---
$casesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
              "WHERE demo.breed IN ['Boston Terrier'] `n" +
              "MATCH (c)<--(diag:diagnosis)`n" +
              "OPTIONAL MATCH (samp:sample)-->(c)`n" +
              "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
              "WITH DISTINCT c, s, demo, diag, co`n" +
              "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
              "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
              "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
              "        coalesce(demo.breed, '') AS Breed ,`n" +
              "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
              "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
              "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
              "        coalesce(demo.sex, '') AS Sex ,`n" +
              "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
              "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
              "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value2 = $casesQuery

# --- Reflect the saved view state: zoomed in to 160%, B2 selected ---
$ws.Activate()
$ws.Range("B2").Select()

$aw = $excel.ActiveWindow
$aw.Zoom = 160
